$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 4707
$ws1.Range("F13").Value = 2866
$ws1.Range("F25").Value = 235

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 4707
$ws4.Range("F14").Value = 2866
$ws4.Range("F26").Value = 235
